$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited range to Text format first so Excel does not
# auto-convert numeric-looking strings (e.g. "65.00") into numbers,
# matching the inline-string/text cells in the original workbook.
$editRange = $ws.Range("B2:E50")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = '58.937.99'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '2.507.82'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '533.11'
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("D6").Value = '135.25'
$ws.Range("E6").Value = '  -2.07%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  -0.64%  '
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("D11").Value = '5.43'
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("D13").Value = '2.954.71'
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("D14").Value = '58.834.28'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").Value = '22.80'
$ws.Range("E15").Value = '  -2.11%  '
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("D17").Value = '2.516.85'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").Value = '322.95'
$ws.Range("E20").Value = '  -0.79%  '
$ws.Range("E21").Value = '  -0.39%  '
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("D23").Value = '65.00'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("E26").Value = '  -1.28%  '
$ws.Range("D27").Value = '7.56'
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("D28").Value = '0.0₃0762'
$ws.Range("E28").Value = '  -2.01%  '
$ws.Range("E29").Value = '  -3.55%  '
$ws.Range("D30").Value = '1.75'
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("D31").Value = '168.83'
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("E33").Value = '  -5.16%  '
$ws.Range("E34").Value = '  -3.74%  '
$ws.Range("E35").Value = '  -1.06%  '
$ws.Range("E36").Value = '  -1.89%  '
$ws.Range("E37").Value = '  -2.90%  '
$ws.Range("D38").Value = '3.56'
$ws.Range("E38").Value = '  -2.23%  '
$ws.Range("E39").Value = '  -4.43%  '
$ws.Range("D40").Value = '281.26'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("E42").Value = '  -0.26%  '
$ws.Range("D43").Value = '5.01'
$ws.Range("E43").Value = '  -5.55%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '129.78'
$ws.Range("E44").Value = '  -0.67%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").Value = '10.92'
$ws.Range("E45").Value = '  +0.57%  '
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("D47").Value = '0.0500'
$ws.Range("E47").Value = '  -2.66%  '
$ws.Range("E48").Value = '  -2.63%  '
$ws.Range("D49").Value = '17.25'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("D50").Value = '1.755.63'
$ws.Range("E50").Value = '  -0.48%  '

# Remove the temporary Text number-format so the cell style matches
# the original (unstyled) cells exactly.
$editRange.ClearFormats()
